# Venue-response analytics: mark which wineries replied with an "x"
# in column C, and a second "x" in column D for the ones that merited
# a follow-up. Row 26 has no reply recorded at all.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$rowsWithC = @(1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,27,28,29,30,31,32)
$rowsWithD = @(5,13,14,17,22,23,30,32)

foreach ($r in $rowsWithC) {
    $ws.Cells.Item($r, 3).Value = "x"
}

foreach ($r in $rowsWithD) {
    $ws.Cells.Item($r, 4).Value = "x"
}

$ws.Range("D23").Select()
